$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the auto "today" date placeholder ("datetimeFigureOut" field)
#    from 19/05/2021 -> 26/05/2021 everywhere it is cached: the slide master
#    and every one of its custom (slide) layouts.
# ---------------------------------------------------------------------------
$oldDate = "19/05/2021"
$newDate = "26/05/2021"

$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the "LLD - Low Level Design" slide that used to sit at position 3
#    (sldId 320). The slide that used to be 4th (sldId 322) shifts up to
#    become the new 3rd slide.
# ---------------------------------------------------------------------------
$p.Slides.Item(3).Delete()
